# Apply the "2e version fichiers" commit changes to the workbook.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Title (row 5) was blank, now carries the resource name.
$meta.Range("B5").Value = "Professionnel"

# Date (row 8) bumped to the new generation timestamp.
$meta.Range("B8").Value = "2025-07-17T14:35:50+00:00"

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 4: Professionnel.typeIdNat-PP / Professionnel.typeIdNat_PP -> Professionnel.typeIdNatPP
$elements.Range("A4").Value = "Professionnel.typeIdNatPP"
$elements.Range("B4").Value = "Professionnel.typeIdNatPP"
$elements.Range("AF4").Value = "Professionnel.typeIdNatPP"

# Row 5: Professionnel.idNat-PS / Professionnel.idNat_PS -> Professionnel.idNatPS
$elements.Range("A5").Value = "Professionnel.idNatPS"
$elements.Range("B5").Value = "Professionnel.idNatPS"
$elements.Range("AF5").Value = "Professionnel.idNatPS"
